$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = '3'
$ws.Range("D2").Value = 0.07915
$ws.Range("E2").Value = -0.08505
$ws.Range("G2").Value = 0.125056709125807
$ws.Range("H2").Value = 0.125056709125807
$ws.Range("I2").Value = 0.04660617693247251
$ws.Range("J2").Value = 0.04331137116494876
$ws.Range("K2").Value = 8.042999999999999
$ws.Range("L2").Value = 0.01403419996510208
$ws.Range("M2").Value = 2.614
$ws.Range("N2").Value = 0.01019103313840156
$ws.Range("O2").Value = 0.3250031082929256
$ws.Range("P2").Value = 2.614
$ws.Range("Q2").Value = 0.01019103313840156
$ws.Range("R2").Value = 0.3250031082929256
$ws.Range("U2").Value = 77.86999999999999
$ws.Range("V2").Value = 0.3035867446393762
$ws.Range("W2").Value = 0.009347653142402546
$ws.Range("X2").Value = 0.08428988055633094
$ws.Range("Y2").Value = -0.07494222741392839
$ws.Range("Z2").Value = 1.430961298377029
$ws.Range("AA2").Value = 0.03752860411899313
$ws.Range("AB2").Value = 0.07928190378314826
$ws.Range("AC2").Value = -0.04175329966415512
$ws.Range("AD2").Value = 67.229
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 67.229
$ws.Range("AG2").Value = -10.64099999999999
$ws.Range("AH2").Value = 0.2076706133834164
$ws.Range("AI2").Value = 0.1485956912576338
$ws.Range("AJ2").Value = -0.04328090490891117
$ws.Range("AK2").Value = -0.02840940946553144
$ws.Range("AL2").Value = 3.128
$ws.Range("AM2").Value = 2.808
$ws.Range("AN2").Value = 1.909914772727273
$ws.Range("AO2").Value = 8.539002557544755
$ws.Range("AP2").Value = -0.3023011363636361
$ws.Range("AQ2").Value = 9.51210826210826

# ---- Row 3 ----
$ws.Range("D3").Value = 0.0563
$ws.Range("E3").Value = -0.09
$ws.Range("G3").Value = 0.12
$ws.Range("H3").Value = 0.12
$ws.Range("I3").Value = 0.2085714285714285
$ws.Range("J3").Value = 0.1643367976043685
$ws.Range("K3").Value = 6.02
$ws.Range("L3").Value = 0.05733333333333333
$ws.Range("M3").Value = 2.51
$ws.Range("N3").Value = 0.06640211640211641
$ws.Range("O3").Value = 0.4169435215946844
$ws.Range("P3").Value = 2.51
$ws.Range("Q3").Value = 0.06640211640211641
$ws.Range("R3").Value = 0.4169435215946844
$ws.Range("U3").Value = 55.9
$ws.Range("V3").Value = 1.478835978835979
$ws.Range("W3").Value = 0.07984084880636604
$ws.Range("X3").Value = 0.07674237486610373
$ws.Range("Y3").Value = 0.003098473940262303
$ws.Range("Z3").Value = 5.303030303030301
$ws.Range("AA3").Value = 0.8714830175989234
$ws.Range("AB3").Value = 0.0762763870143294
$ws.Range("AC3").Value = 0.7952066305845941
$ws.Range("AD3").Value = 0.829
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0.829
$ws.Range("AG3").Value = -55.071
$ws.Range("AH3").Value = 0.02146056071863108
$ws.Range("AI3").Value = 0.01041077999221389
$ws.Range("AJ3").Value = 3.188639916623241
$ws.Range("AK3").Value = -2.320831050613174
$ws.Range("AL3").Value = 0.128
$ws.Range("AM3").Value = 0.128
$ws.Range("AN3").Value = 0.03700892857142857
$ws.Range("AO3").Value = 171.09375
$ws.Range("AP3").Value = -2.458526785714286
$ws.Range("AQ3").Value = 171.09375

# ---- Row 4 ----
$ws.Range("B4").Value = 'Britam Holdings Plc (NASE:BRIT)'
$ws.Range("D4").Value = 0.102
$ws.Range("E4").Value = -0.0801
$ws.Range("G4").Value = 0.1864352298296368
$ws.Range("H4").Value = 0.1864352298296368
$ws.Range("I4").Value = 0.03162970106075217
$ws.Range("J4").Value = 0.03162970106075217
$ws.Range("K4").Value = 2.35
$ws.Range("L4").Value = 0.007553841208614593
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("U4").Value = 18.4
$ws.Range("V4").Value = 0.1097197376267144
$ws.Range("W4").Value = 0.009347653142402546
$ws.Range("X4").Value = 0.08428988055633094
$ws.Range("Y4").Value = -0.07494222741392839
$ws.Range("Z4").Value = 1.186498855835241
$ws.Range("AA4").Value = 0.03752860411899313
$ws.Range("AB4").Value = 0.07928190378314826
$ws.Range("AC4").Value = -0.04175329966415512
$ws.Range("AD4").Value = 30.9
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 30.9
$ws.Range("AG4").Value = 12.5
$ws.Range("AH4").Value = 0.1555891238670695
$ws.Range("AI4").Value = 0.1158170914542729
$ws.Range("AJ4").Value = 0.06936736958934518
$ws.Range("AK4").Value = 0.05032206119162641
$ws.Range("AL4").Value = 3
$ws.Range("AM4").Value = 3
$ws.Range("AN4").Value = 2.4140625
$ws.Range("AO4").Value = 3.28
$ws.Range("AP4").Value = 0.9765625
$ws.Range("AQ4").Value = 3.28
$ws.Range("T4").ClearContents()

# ---- Row 5 ----
$ws.Range("G5").Value = 0.006815286624203822
$ws.Range("H5").Value = 0.006815286624203822
$ws.Range("I5").Value = -0.03203821656050956
$ws.Range("J5").Value = -0.03203821656050956
$ws.Range("K5").Value = -0.327
$ws.Range("L5").Value = -0.002082802547770701
$ws.Range("M5").Value = 0.104
$ws.Range("N5").Value = 0.00203921568627451
$ws.Range("O5").Value = -0.3180428134556575
$ws.Range("P5").Value = 0.104
$ws.Range("Q5").Value = 0.00203921568627451
$ws.Range("R5").Value = -0.3180428134556575
$ws.Range("U5").Value = 3.57
$ws.Range("V5").Value = 0.06999999999999999
$ws.Range("W5").Value = -0.004336870026525199
$ws.Range("X5").Value = 0.10808743258955
$ws.Range("Y5").Value = -0.1124243026160752
$ws.Range("Z5").Value = 1.324894514767933
$ws.Range("AA5").Value = -0.04244725738396625
$ws.Range("AB5").Value = 0.08321496599142086
$ws.Range("AC5").Value = -0.1256622233753871
$ws.Range("AD5").Value = 35.5
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 35.5
$ws.Range("AG5").Value = 31.93
$ws.Range("AH5").Value = 0.4104046242774567
$ws.Range("AI5").Value = 0.3349056603773585
$ws.Range("AJ5").Value = 0.385023513806825
$ws.Range("AK5").Value = 0.3117250805428097
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = -0.32
$ws.Range("AQ5").Value = 15.71875
$ws.Range("AN5").ClearContents()
$ws.Range("AO5").ClearContents()
$ws.Range("AP5").ClearContents()

# ---- Remove row 6 entirely ----
$ws.Rows(6).Delete()
